$d = $word.ActiveDocument

$pairs = @(
    @("317÷8=39, 5", "822÷7=117, 3"),
    @("269÷2=134, 1", "815÷4=203, 3"),
    @("710÷5=142, 0", "712÷5=142, 2"),
    @("253÷9=28, 1", "740÷4=185, 0"),
    @("886÷3=295, 1", "590÷2=295, 0"),
    @("815÷5=163, 0", "635÷6=105, 5"),
    @("513÷8=64, 1", "485÷2=242, 1"),
    @("707÷9=78, 5", "904÷4=226, 0"),
    @("841÷8=105, 1", "622÷2=311, 0"),
    @("603÷2=301, 1", "688÷2=344, 0"),
    @("156÷9=17, 3", "261÷2=130, 1"),
    @("929÷3=309, 2", "414÷2=207, 0"),
    @("809÷3=269, 2", "131÷2=65, 1"),
    @("231÷3=77, 0", "433÷5=86, 3"),
    @("468÷6=78, 0", "377÷2=188, 1"),
    @("101÷2=50, 1", "524÷8=65, 4"),
    @("517÷6=86, 1", "431÷5=86, 1"),
    @("465÷4=116, 1", "604÷3=201, 1"),
    @("648÷3=216, 0", "573÷5=114, 3"),
    @("530÷8=66, 2", "225÷7=32, 1"),
    @("304÷8=38, 0", "323÷5=64, 3"),
    @("765÷8=95, 5", "265÷2=132, 1"),
    @("994÷9=110, 4", "262÷5=52, 2"),
    @("467÷4=116, 3", "453÷8=56, 5"),
    @("786÷4=196, 2", "773÷7=110, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
